# Insert two new rows for a new weekly data point (Brocoli, Terminal La
# Palmera de La Serena). This pushes the existing rows 245-372 down to
# 247-374 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("245:246").Insert()

# New row 245 (Primera)
$ws.Cells.Item(245, 1).Value = 8
$ws.Cells.Item(245, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(245, 3).Value = "Coquimbo"
$ws.Cells.Item(245, 4).Value = 44452
$ws.Cells.Item(245, 5).Value = 4
$ws.Cells.Item(245, 6).Value = 100112023
$ws.Cells.Item(245, 7).Value = "Brócoli"
$ws.Cells.Item(245, 8).Value = "Sin especificar"
$ws.Cells.Item(245, 9).Value = "Primera"
$ws.Cells.Item(245, 10).Value = 3000
$ws.Cells.Item(245, 11).Value = 600
$ws.Cells.Item(245, 12).Value = 700
$ws.Cells.Item(245, 13).Value = 650
$ws.Cells.Item(245, 14).Value = "`$/unidad"
$ws.Cells.Item(245, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(245, 16).Value = 650
$ws.Cells.Item(245, 17).Value = 1
$ws.Cells.Item(245, 18).Value = "Hortaliza"

# New row 246 (Segunda)
$ws.Cells.Item(246, 1).Value = 8
$ws.Cells.Item(246, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(246, 3).Value = "Coquimbo"
$ws.Cells.Item(246, 4).Value = 44452
$ws.Cells.Item(246, 5).Value = 4
$ws.Cells.Item(246, 6).Value = 100112023
$ws.Cells.Item(246, 7).Value = "Brócoli"
$ws.Cells.Item(246, 8).Value = "Sin especificar"
$ws.Cells.Item(246, 9).Value = "Segunda"
$ws.Cells.Item(246, 10).Value = 1400
$ws.Cells.Item(246, 11).Value = 500
$ws.Cells.Item(246, 12).Value = 550
$ws.Cells.Item(246, 13).Value = 525
$ws.Cells.Item(246, 14).Value = "`$/unidad"
$ws.Cells.Item(246, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(246, 16).Value = 525
$ws.Cells.Item(246, 17).Value = 1
$ws.Cells.Item(246, 18).Value = "Hortaliza"
